# Fix upload-users worksheet: replace sample data rows with the
# George McFly / Marty McFly / Doc Brown records and drop the old
# third data row (former row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("I2").Value = "'0010010"
$ws.Range("B2").Value = "George McFly"
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = "Doc Brown"
$ws.Range("E2").Value = "Chemistry"
$ws.Range("F2").Value = "ArtSci"
$ws.Range("G2").Value = "Queen's"

# --- Row 3 -------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Marty McFly"
$ws.Range("C3").Value = 14
$ws.Range("D3").Value = "Doc Brown"
$ws.Range("E3").Value = "Chemistry"
$ws.Range("F3").Value = "ArtSci"
$ws.Range("G3").Value = "Queen's"
$ws.Range("I3").Value = "'0100011"
$ws.Range("H3").Value = "'50"

# H2 filled in after the rest of row 2/3, matching original edit order
$ws.Range("H2").Value = "Academic Machine Dependent"

# --- Remove the old 4th data row ----------------------------------------
$ws.Rows("4:4").Delete()

# --- Column sizing -------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.109375
$ws.Columns.Item(8).ColumnWidth = 35.109375

# --- Selection -------------------------------------------------------
$null = $ws.Range("H3").Select()

Write-Output "done"
